# Scheduled market-board data refresh for Chocobo_Profits workbook.
# Source: commit "chore: update Sheets via scheduled runner".
# Refreshes currentAveragePrice(NQ/HQ) and LevePrice/LeveProfit columns (H:N)
# for the leves whose market data moved since the last pull.

$wb = $excel.ActiveWorkbook

# ALC row 4: Root Rush
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 141.2
$ws.Range("I4").Value = 112.44444
$ws.Range("K4").Value = 112.44444
$ws.Range("M4").Value = 1.55556

# ALC row 5: Met a Sticky End
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 433.66666
$ws.Range("I5").Value = 101
$ws.Range("K5").Value = 101
$ws.Range("M5").Value = 14

# ALC row 123: Nearly Bare
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 43780
$ws.Range("J123").Value = 43780
$ws.Range("L123").Value = 43780
$ws.Range("N123").Value = -53580

# ARM row 61: Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1307.762
$ws.Range("I61").Value = 1168.15
$ws.Range("J61").Value = 4100
$ws.Range("K61").Value = 1168.15
$ws.Range("L61").Value = 4100
$ws.Range("M61").Value = -956.1500000000001
$ws.Range("N61").Value = -4524

# ARM row 63: Rivets Run through It
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8659877
$ws.Range("I63").Value = 13852993
$ws.Range("K63").Value = 13852993
$ws.Range("M63").Value = -13852307

# ARM row 66: A Riveting Revival (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 8659877
$ws.Range("I66").Value = 13852993
$ws.Range("K66").Value = 69264965
$ws.Range("M66").Value = -69261533

# ARM row 136: Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1307.762
$ws.Range("I136").Value = 1168.15
$ws.Range("J136").Value = 4100
$ws.Range("K136").Value = 3504.45
$ws.Range("L136").Value = 12300
$ws.Range("M136").Value = -954.4500000000003
$ws.Range("N136").Value = -17400

# BSM row 10: Bring Me the Head Knife of Al'bedo Derssia
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 4366
$ws.Range("J10").Value = 4366
$ws.Range("L10").Value = 4366
$ws.Range("N10").Value = -4646

# CRP row 11: Leaving without Leave
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 24500
$ws.Range("J11").Value = 24500
$ws.Range("L11").Value = 24500
$ws.Range("N11").Value = -24780

# CRP row 48: The Cold, Cold Ground
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 45005
$ws.Range("J48").Value = 45005
$ws.Range("L48").Value = 45005
$ws.Range("N48").Value = -45957

# CUL row 14: Keep Your Powder Dry
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 566.6957
$ws.Range("I14").Value = 566.6957
$ws.Range("K14").Value = 1700.0871
$ws.Range("M14").Value = -1527.0871

# CUL row 34: Fever Pitch
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 11849.315
$ws.Range("I34").Value = 27587.5
$ws.Range("J34").Value = 7652.467
$ws.Range("K34").Value = 82762.5
$ws.Range("L34").Value = 22957.401
$ws.Range("M34").Value = -82678.5
$ws.Range("N34").Value = -23125.401

# CUL row 39: Bloody Good Tart, This
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 9929.637000000001
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 9929.637000000001
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 29788.911
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -30376.911

# CUL row 55: Pagan Pastries
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 4746.0415
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 4865.4346
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 14596.3038
$ws.Range("M55").Value = -5823
$ws.Range("N55").Value = -14950.3038

# CUL row 98: Sweet Kiss of Death
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 266
$ws.Range("I98").Value = 226.5
$ws.Range("J98").Value = 272.58334
$ws.Range("K98").Value = 679.5
$ws.Range("L98").Value = 817.7500200000001
$ws.Range("M98").Value = 818.5
$ws.Range("N98").Value = -3813.75002

# GSM row 2: Copper and Robbers
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 74.44444
$ws.Range("I2").Value = 66.666664
$ws.Range("K2").Value = 66.666664
$ws.Range("M2").Value = 46.333336

# GSM row 5: Hora at Me
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10945
$ws.Range("J5").Value = 10986
$ws.Range("L5").Value = 10986
$ws.Range("N5").Value = -11210

# GSM row 12: Horn of Plenty
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 20000
$ws.Range("J12").Value = 20000
$ws.Range("L12").Value = 20000
$ws.Range("N12").Value = -20280

# GSM row 14: All That Glitters
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4310050.5
$ws.Range("I14").Value = 5375063.5
$ws.Range("J14").Value = 49999.5
$ws.Range("K14").Value = 5375063.5
$ws.Range("L14").Value = 49999.5
$ws.Range("M14").Value = -5374895.5
$ws.Range("N14").Value = -50335.5

# GSM row 74: The Unfortunate Retirony
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 32000
$ws.Range("J74").Value = 32000
$ws.Range("L74").Value = 32000
$ws.Range("N74").Value = -33872

# GSM row 77: Life Ends at Retirement (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 32000
$ws.Range("J77").Value = 32000
$ws.Range("L77").Value = 96000
$ws.Range("N77").Value = -105360

# GSM row 126: Gold Rush Order
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2985.12
$ws.Range("I126").Value = 2970.6526
$ws.Range("J126").Value = 3260
$ws.Range("K126").Value = 8911.9578
$ws.Range("L126").Value = 9780
$ws.Range("M126").Value = -6441.9578
$ws.Range("N126").Value = -14720

# GSM row 132: On Board for Lar
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3237.7942
$ws.Range("I132").Value = 2135.8
$ws.Range("J132").Value = 4107.7896
$ws.Range("K132").Value = 6407.400000000001
$ws.Range("L132").Value = 12323.3688
$ws.Range("M132").Value = -3877.400000000001
$ws.Range("N132").Value = -17383.3688

# LTW row 122: Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6460.533
$ws.Range("I122").Value = 4129.7144
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 12389.1432
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -9939.143199999999
$ws.Range("N122").Value = -30400

# WVR row 2: The Unmentionables
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 16445
$ws.Range("I2").Value = 1002
$ws.Range("K2").Value = 1002
$ws.Range("M2").Value = -890

# WVR row 4: Not Cool Enough
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2833.6667
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 3200.4
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 3200.4
$ws.Range("M4").Value = -887
$ws.Range("N4").Value = -3426.4

# WVR row 5: Hire in the Blood
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 125630700
$ws.Range("I5").Value = 333338340
$ws.Range("J5").Value = 1006119.8
$ws.Range("K5").Value = 333338340
$ws.Range("L5").Value = 1006119.8
$ws.Range("M5").Value = -333338228
$ws.Range("N5").Value = -1006343.8

# WVR row 10: Just for Kecks
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 5000
$ws.Range("J10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("N10").Value = -5338

# WVR row 50: Cool to Be Southern
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 29000
$ws.Range("J50").Value = 29000
$ws.Range("L50").Value = 29000
$ws.Range("N50").Value = -30262

# WVR row 132: Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13891071
$ws.Range("I132").Value = 1169.5
$ws.Range("J132").Value = 33336934
$ws.Range("K132").Value = 3508.5
$ws.Range("L132").Value = 100010802
$ws.Range("M132").Value = -978.5
$ws.Range("N132").Value = -100015862
